$wb = $excel.ActiveWorkbook

$oldGuid = "416eb773-b388-42a3-9539-864b55141604"
$newGuid = "899ee086-8242-4535-95c3-0bab5ea32bdf"

$newHash = "ee505f28e96fef06e3947d82f0999e12df38f31f"

$hyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58273e50f01dd03a723ef7c1d1c263c7bf2124ee/e2e/$oldGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-15 22:53:53"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddr, "", "", "e2e\$newGuid.md")

# ---- zh-cn sheet ----
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-15 22:53:48"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddr, "", "", "$newGuid.md")

# ---- de-de sheet ----
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-15 22:53:53"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddr, "", "", "$newGuid.md")
